$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Build the new "OpenMP" results section starting at row 12
#    Doing this BEFORE touching F1 means this text becomes a brand new shared
#    string entry (appended at the end), while F1's old string slot is freed
#    up to be renamed below - matching how the original authors' diff shows
#    the shared-string table evolving.

# Section title row, merged A12:D12, centered + wrapped text
$title = $ws.Range("A12:D12")
$title.Merge()
$title.HorizontalAlignment = -4108   # xlCenter
$title.WrapText = $true
$ws.Range("A12").Value = "OpenMP, primes until 1750, upper limit: 2^1000000"

# 2) Update the F1 header text: "C++, primes until 1750..." -> "C++, primes until 2000..."
#    (F1 was the sole reference to the old string slot, so its text is replaced in place)
$ws.Range("F1").Value = "C++, primes until 2000, upper limit: 2^1000000"

# Header row 13 (same headers as row 2/13)
$ws.Range("A13").Value = "Processes"
$ws.Range("B13").Value = "Execution time (s)"
$ws.Range("C13").Value = "Speedup"
$ws.Range("D13").Value = "Efficiency"

# Data rows 14-21
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 194.28

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 202.59
$ws.Range("C15").Formula = "=B14/B15"
$ws.Range("D15").Formula = "=C15/A15"

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 217.84
$ws.Range("C16").Formula = "=B14/B16"
$ws.Range("D16").Formula = "=C16/A16"

$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 231.65
$ws.Range("C17").Formula = "=B14/B17"
$ws.Range("D17").Formula = "=C17/A17"

$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 247.35
$ws.Range("C18").Formula = "=B14/B18"
$ws.Range("D18").Formula = "=C18/A18"

$ws.Range("A19").Value = 6
$ws.Range("B19").Value = 281.20999999999998
$ws.Range("C19").Formula = "=B14/B19"
$ws.Range("D19").Formula = "=C19/A19"

$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 298.5
$ws.Range("C20").Formula = "=B14/B20"
$ws.Range("D20").Formula = "=C20/A20"

$ws.Range("A21").Value = 8
$ws.Range("B21").Value = 305.33999999999997
$ws.Range("C21").Formula = "=B14/B21"
$ws.Range("D21").Formula = "=C21/A21"

# 3) Fix up the view: selection on F1:I1 like the header above
$ws.Range("F1:I1").Select()
